$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# C1: give it a "middle of header box" border (top + bottom only)
$c1 = $ws1.Range("C1")
$c1.ClearFormats()
$c1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$c1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# D1: give it the "right edge of header box" border (right + top + bottom)
$d1 = $ws1.Range("D1")
$d1.ClearFormats()
$d1.Borders.Item(10).LineStyle = 1  # xlEdgeRight
$d1.Borders.Item(8).LineStyle = 1   # xlEdgeTop
$d1.Borders.Item(9).LineStyle = 1   # xlEdgeBottom

# Rename "fedcore" header to "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

$c1b = $ws2.Range("C1")
$c1b.ClearFormats()
$c1b.Borders.Item(8).LineStyle = 1
$c1b.Borders.Item(9).LineStyle = 1

$d1b = $ws2.Range("D1")
$d1b.ClearFormats()
$d1b.Borders.Item(10).LineStyle = 1
$d1b.Borders.Item(8).LineStyle = 1
$d1b.Borders.Item(9).LineStyle = 1

$f1b = $ws2.Range("F1")
$f1b.ClearFormats()
$f1b.Borders.Item(8).LineStyle = 1
$f1b.Borders.Item(9).LineStyle = 1

$g1b = $ws2.Range("G1")
$g1b.ClearFormats()
$g1b.Borders.Item(10).LineStyle = 1
$g1b.Borders.Item(8).LineStyle = 1
$g1b.Borders.Item(9).LineStyle = 1

# Rename "fedcore" headers to "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
